# Auto-applies the Linea 141 horarios scrape-refresh update (commit: "Horarios actualizados Linea 141 - 1153").
# New scrape snapshot (08:04:39) merged into each sheet's time-sorted table; header/meta rows refreshed,
# a couple of late-breaking rows from the previous snapshot (07:28:23) got re-ordered, and newly scraped
# rows were appended so every sheet stays sorted by Hora_Llegada (column B).
$wb = $excel.ActiveWorkbook

# --- Sheet "LP1912" ---
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(2,1).Value = "Última actualización: 08:04:39"
$ws.Cells.Item(3,1).Value = "Total filas: 73"
$ws.Cells.Item(22,3).Value = "215C_EL PATO"
$ws.Cells.Item(23,3).Value = "14_ABASTO"
$ws.Cells.Item(48,1).Value = "08:04:39"
$ws.Cells.Item(48,2).Value = "08:04"
$ws.Cells.Item(48,3).Value = "16_SANTA ANA"
$ws.Cells.Item(48,4).Value = 0
$ws.Cells.Item(49,1).Value = "07:28:23"
$ws.Cells.Item(49,2).Value = "08:08"
$ws.Cells.Item(49,4).Value = 40
$ws.Cells.Item(50,1).Value = "06:55:48"
$ws.Cells.Item(50,2).Value = "08:09"
$ws.Cells.Item(50,3).Value = "23_HERNANDEZ"
$ws.Cells.Item(50,4).Value = 74
$ws.Cells.Item(51,1).Value = "08:04:39"
$ws.Cells.Item(51,2).Value = "08:10"
$ws.Cells.Item(51,3).Value = "16_SANTA ANA"
$ws.Cells.Item(51,4).Value = 6
$ws.Cells.Item(52,1).Value = "08:04:39"
$ws.Cells.Item(52,2).Value = "08:12"
$ws.Cells.Item(52,3).Value = "15_ABASTO"
$ws.Cells.Item(52,4).Value = 8
$ws.Cells.Item(53,1).Value = "08:04:39"
$ws.Cells.Item(53,2).Value = "08:13"
$ws.Cells.Item(53,3).Value = "10_OLMOS"
$ws.Cells.Item(53,4).Value = 9
$ws.Cells.Item(54,1).Value = "08:04:39"
$ws.Cells.Item(54,2).Value = "08:21"
$ws.Cells.Item(54,3).Value = "26_HERNANDEZ"
$ws.Cells.Item(54,4).Value = 17
$ws.Cells.Item(55,1).Value = "08:04:39"
$ws.Cells.Item(55,2).Value = "08:22"
$ws.Cells.Item(55,3).Value = "16_P MOR-SANTA ANA"
$ws.Cells.Item(55,4).Value = 18
$ws.Cells.Item(56,1).Value = "08:04:39"
$ws.Cells.Item(56,2).Value = "08:23"
$ws.Cells.Item(56,3).Value = "215B_EL PATO"
$ws.Cells.Item(56,4).Value = 19
$ws.Cells.Item(57,1).Value = "08:04:39"
$ws.Cells.Item(57,2).Value = "08:27"
$ws.Cells.Item(57,3).Value = "84_COLONIA URQUIZA-ESC 49"
$ws.Cells.Item(57,4).Value = 23
$ws.Cells.Item(58,1).Value = "08:04:39"
$ws.Cells.Item(58,2).Value = "08:33"
$ws.Cells.Item(58,3).Value = "10_OLMOS"
$ws.Cells.Item(58,4).Value = 29
$ws.Cells.Item(59,1).Value = "08:04:39"
$ws.Cells.Item(59,2).Value = "08:34"
$ws.Cells.Item(59,3).Value = "23_HERNANDEZ"
$ws.Cells.Item(59,4).Value = 30
$ws.Cells.Item(60,1).Value = "06:55:48"
$ws.Cells.Item(60,2).Value = "08:42"
$ws.Cells.Item(60,3).Value = "81_EL PELIGRO"
$ws.Cells.Item(60,4).Value = 107
$ws.Cells.Item(61,1).Value = "08:04:39"
$ws.Cells.Item(61,2).Value = "08:43"
$ws.Cells.Item(61,3).Value = "14_ABASTO"
$ws.Cells.Item(61,4).Value = 39
$ws.Cells.Item(62,2).Value = "08:50"
$ws.Cells.Item(62,3).Value = "81_EL PELIGRO"
$ws.Cells.Item(62,4).Value = 82
$ws.Cells.Item(63,1).Value = "08:04:39"
$ws.Cells.Item(63,2).Value = "08:54"
$ws.Cells.Item(63,3).Value = "17_ROMERO"
$ws.Cells.Item(63,4).Value = 50
$ws.Cells.Item(64,1).Value = "08:04:39"
$ws.Cells.Item(64,2).Value = "09:01"
$ws.Cells.Item(64,3).Value = "215A_EL PATO"
$ws.Cells.Item(64,4).Value = 57
$ws.Cells.Item(65,1).Value = "08:04:39"
$ws.Cells.Item(65,2).Value = "09:02"
$ws.Cells.Item(65,3).Value = "23_HERNANDEZ"
$ws.Cells.Item(65,4).Value = 58
$ws.Cells.Item(65,5).Value = "LP1912"
$ws.Cells.Item(66,1).Value = "08:04:39"
$ws.Cells.Item(66,2).Value = "09:03"
$ws.Cells.Item(66,3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(66,4).Value = 59
$ws.Cells.Item(66,5).Value = "LP1912"
$ws.Cells.Item(67,1).Value = "08:04:39"
$ws.Cells.Item(67,2).Value = "09:10"
$ws.Cells.Item(67,3).Value = "16_P MOR-SANTA ANA"
$ws.Cells.Item(67,4).Value = 66
$ws.Cells.Item(67,5).Value = "LP1912"
$ws.Cells.Item(68,1).Value = "08:04:39"
$ws.Cells.Item(68,2).Value = "09:16"
$ws.Cells.Item(68,3).Value = "27_EL RETIRO"
$ws.Cells.Item(68,4).Value = 72
$ws.Cells.Item(68,5).Value = "LP1912"
$ws.Cells.Item(69,1).Value = "07:28:23"
$ws.Cells.Item(69,2).Value = "09:17"
$ws.Cells.Item(69,3).Value = "27_EL RETIRO"
$ws.Cells.Item(69,4).Value = 109
$ws.Cells.Item(69,5).Value = "LP1912"
$ws.Cells.Item(70,1).Value = "08:04:39"
$ws.Cells.Item(70,2).Value = "09:21"
$ws.Cells.Item(70,3).Value = "26_HERNANDEZ"
$ws.Cells.Item(70,4).Value = 77
$ws.Cells.Item(70,5).Value = "LP1912"
$ws.Cells.Item(71,1).Value = "08:04:39"
$ws.Cells.Item(71,2).Value = "09:22"
$ws.Cells.Item(71,3).Value = "17_ROMERO"
$ws.Cells.Item(71,4).Value = 78
$ws.Cells.Item(71,5).Value = "LP1912"
$ws.Cells.Item(72,1).Value = "07:28:23"
$ws.Cells.Item(72,2).Value = "09:23"
$ws.Cells.Item(72,3).Value = "17_ROMERO"
$ws.Cells.Item(72,4).Value = 115
$ws.Cells.Item(72,5).Value = "LP1912"
$ws.Cells.Item(73,1).Value = "08:04:39"
$ws.Cells.Item(73,2).Value = "09:23"
$ws.Cells.Item(73,3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(73,4).Value = 79
$ws.Cells.Item(73,5).Value = "LP1912"
$ws.Cells.Item(74,1).Value = "08:04:39"
$ws.Cells.Item(74,2).Value = "09:25"
$ws.Cells.Item(74,3).Value = "81_EL PELIGRO"
$ws.Cells.Item(74,4).Value = 81
$ws.Cells.Item(74,5).Value = "LP1912"
$ws.Cells.Item(75,1).Value = "08:04:39"
$ws.Cells.Item(75,2).Value = "09:32"
$ws.Cells.Item(75,3).Value = "15_ABASTO"
$ws.Cells.Item(75,4).Value = 88
$ws.Cells.Item(75,5).Value = "LP1912"
$ws.Cells.Item(76,1).Value = "08:04:39"
$ws.Cells.Item(76,2).Value = "09:33"
$ws.Cells.Item(76,3).Value = "10_OLMOS"
$ws.Cells.Item(76,4).Value = 89
$ws.Cells.Item(76,5).Value = "LP1912"
$ws.Cells.Item(77,1).Value = "08:04:39"
$ws.Cells.Item(77,2).Value = "09:41"
$ws.Cells.Item(77,3).Value = "215C_EL PATO"
$ws.Cells.Item(77,4).Value = 97
$ws.Cells.Item(77,5).Value = "LP1912"
$ws.Cells.Item(78,1).Value = "08:04:39"
$ws.Cells.Item(78,2).Value = "09:43"
$ws.Cells.Item(78,3).Value = "14_ABASTO"
$ws.Cells.Item(78,4).Value = 99
$ws.Cells.Item(78,5).Value = "LP1912"

# --- Sheet "LP1912-215" ---
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(2,1).Value = "Última actualización: 08:04:39"
$ws.Cells.Item(3,1).Value = "Total filas: 10"
$ws.Cells.Item(13,1).Value = "08:04:39"
$ws.Cells.Item(13,4).Value = 19
$ws.Cells.Item(14,1).Value = "08:04:39"
$ws.Cells.Item(14,4).Value = 57
$ws.Cells.Item(15,1).Value = "08:04:39"
$ws.Cells.Item(15,2).Value = "09:41"
$ws.Cells.Item(15,3).Value = "215C_EL PATO"
$ws.Cells.Item(15,4).Value = 97
$ws.Cells.Item(15,5).Value = "LP1912"

# --- Sheet "6203-6173" ---
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(2,1).Value = "Última actualización: 08:04:39"
$ws.Cells.Item(3,1).Value = "Total filas: 17"
$ws.Cells.Item(17,1).Value = "08:04:39"
$ws.Cells.Item(17,2).Value = "08:25"
$ws.Cells.Item(17,3).Value = "215C_LA PLATA"
$ws.Cells.Item(17,4).Value = 21
$ws.Cells.Item(17,5).Value = "L6203"
$ws.Cells.Item(18,1).Value = "06:55:48"
$ws.Cells.Item(18,2).Value = "08:31"
$ws.Cells.Item(18,4).Value = 96
$ws.Cells.Item(19,1).Value = "08:04:39"
$ws.Cells.Item(19,2).Value = "08:35"
$ws.Cells.Item(19,3).Value = "215A_LA PLATA"
$ws.Cells.Item(19,4).Value = 31
$ws.Cells.Item(19,5).Value = "L6173"
$ws.Cells.Item(20,1).Value = "08:04:39"
$ws.Cells.Item(20,2).Value = "09:08"
$ws.Cells.Item(20,3).Value = "215D_LA PLATA"
$ws.Cells.Item(20,4).Value = 64
$ws.Cells.Item(20,5).Value = "L6203"
$ws.Cells.Item(21,1).Value = "07:28:23"
$ws.Cells.Item(21,2).Value = "09:09"
$ws.Cells.Item(21,3).Value = "215D_LA PLATA"
$ws.Cells.Item(21,4).Value = 101
$ws.Cells.Item(21,5).Value = "L6203"
$ws.Cells.Item(22,1).Value = "08:04:39"
$ws.Cells.Item(22,2).Value = "10:02"
$ws.Cells.Item(22,3).Value = "215B_LP-P MOR-40 Y 115"
$ws.Cells.Item(22,4).Value = 118
$ws.Cells.Item(22,5).Value = "L6173"
